$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Password" column: header + the two teachers' passwords
$ws.Range("D1").Value = "Password"
$ws.Range("D2").Value = "rks_1987"
$ws.Range("D3").Value = "prk_1985"

# The Gmail ID column now needs to fit its (already long) contents alongside
# the new column, so auto-fit its width like the author did
$ws.Columns.Item(3).AutoFit() | Out-Null

# Leave the workbook scrolled/selected on the new column, matching where the
# author's cursor ended up when they saved
$ws.Range("D17").Select()
